$d = $word.ActiveDocument

# 1) Replace the employee name
#    " Milton Stiven Gonzalez Pinzon" -> " KELLYN JOHANNA DELGADO JAIMES"
$d.Content.Find.Execute(" Milton Stiven Gonzalez Pinzon", $true, $false, $false, $false, $false, $true, 1, $false, " KELLYN JOHANNA DELGADO JAIMES", 2) | Out-Null

# 2) Replace the Cedula / ID number
#    1013099140 -> 1127586868
$d.Content.Find.Execute("1013099140", $true, $false, $false, $false, $false, $true, 1, $false, "1127586868", 2) | Out-Null

# 3) Replace the job title and fill in the (previously empty) contract-type
#    sentence that follows it.
#    "Director de Mantenimiento" -> "ASISTENTE"
#    <empty run> -> "Mediante un contrato a Término Fijo."
$d.Content.Find.Execute("Director de Mantenimiento. ", $true, $false, $false, $false, $false, $true, 1, $false, "ASISTENTE. Mediante un contrato a Término Fijo.", 2) | Out-Null

# 4) Fill in the (previously empty) contract dates + salary sentence.
#    <empty run> -> "Desde el 2024-02-15 hasta el ."
#    " " -> " devengando un salario de $ 123213."
#    These two runs sit right before the paragraph mark, after a distinctly
#    formatted " " run (w:spacing=16) that must stay untouched, so address
#    them positionally (by character offset) rather than via Find/Replace.
$p = $d.Paragraphs.Item(10)
$paraEnd = $p.Range.End
$tailRun = $d.Range($paraEnd - 2, $paraEnd - 1)
$tailRun.Text = "Desde el 2024-02-15 hasta el . devengando un salario de `$ 123213."

# 5) Replace the day-of-month in the closing date line
#    "(09) días del mes de (febrero)" -> "(16) días del mes de (febrero)"
$d.Content.Find.Execute("(09) días del mes de (febrero)", $true, $false, $false, $false, $false, $true, 1, $false, "(16) días del mes de (febrero)", 2) | Out-Null
